$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to D and E columns first so numeric-looking
# strings (e.g. "1.000", "30.626.61") are stored as text, matching the
# original inline-string cell content instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.592.48"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "1.923.19"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "247.32"
$ws.Range("E5").Value = "  +2.80%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "0.4732"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").Value = "0.2915"
$ws.Range("E8").Value = "  +1.14%  "

$ws.Range("D9").Value = "0.06823"
$ws.Range("E9").Value = "  +2.75%  "

$ws.Range("D10").Value = "105.42"
$ws.Range("E10").Value = "  -2.20%  "

$ws.Range("D11").Value = "18.38"
$ws.Range("E11").Value = "  -4.12%  "

$ws.Range("D12").Value = "1.921.64"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").Value = "0.07733"
$ws.Range("E13").Value = "  +1.57%  "

$ws.Range("D14").Value = "5.321"
$ws.Range("E14").Value = "  +3.06%  "

$ws.Range("D15").Value = "0.6708"
$ws.Range("E15").Value = "  +1.19%  "

$ws.Range("D16").Value = "292.75"
$ws.Range("E16").Value = "  -5.20%  "

$ws.Range("D17").Value = "30.626.61"

$ws.Range("D18").Value = "0.000007636"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.548"
$ws.Range("E21").Value = "  +4.65%  "

$ws.Range("D22").Value = "2.178.69"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "6.466"
$ws.Range("E24").Value = "  +2.55%  "

$ws.Range("D25").Value = "9.518"
$ws.Range("E25").Value = "  +2.02%  "

$ws.Range("D26").Value = "167.64"
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").Value = "20.86"
$ws.Range("E27").Value = "  +2.78%  "

$ws.Range("D28").Value = "2.133"
$ws.Range("E28").Value = "  +3.85%  "

$ws.Range("D29").Value = "0.1070"
$ws.Range("E29").Value = "  -3.94%  "

$ws.Range("D30").Value = "1.403"
$ws.Range("E30").Value = "  +2.16%  "

$ws.Range("D31").Value = "4.200"
$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("D32").Value = "4.055"
$ws.Range("E32").Value = "  +2.86%  "

$ws.Range("D33").Value = "0.05035"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "0.7346"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").Value = "1.146"
$ws.Range("E35").Value = "  -1.08%  "

$ws.Range("D36").Value = "0.02061"
$ws.Range("E36").Value = "  +4.81%  "

$ws.Range("D37").Value = "0.9993"

$ws.Range("D38").Value = "2.736"
$ws.Range("E38").Value = "  -0.85%  "

$ws.Range("E39").Value = "  -0.41%  "

$ws.Range("D40").Value = "111.84"
$ws.Range("E40").Value = "  +3.81%  "

$ws.Range("D41").Value = "2.044"
$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("D42").Value = "0.4448"
$ws.Range("E42").Value = "  +6.00%  "

$ws.Range("D43").Value = "0.8711"
$ws.Range("E43").Value = "  -1.14%  "

$ws.Range("D44").Value = "5.893"
$ws.Range("E44").Value = "  +1.22%  "

$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").Value = "67.85"
$ws.Range("E46").Value = "  -3.94%  "

$ws.Range("D47").Value = "7.295"
$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").Value = "9.425"
$ws.Range("E48").Value = "  +1.68%  "

$ws.Range("D49").Value = "0.1251"
$ws.Range("E49").Value = "  +2.87%  "

$ws.Range("D50").Value = "47.73"
$ws.Range("E50").Value = "  +12.30%  "

$ws.Range("D51").Value = "35.19"
$ws.Range("E51").Value = "  +0.74%  "
